$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 10: fold the "L..Q" (old matmult-style) columns into the
# "D..J" (sum/mult style) columns used by the rest of the sheet, and
# add a new multiplier column I10.
# ---------------------------------------------------------------------
$ws.Cells.Item(10,4).Value = 2047
$ws.Cells.Item(10,5).Value = 2050
$ws.Cells.Item(10,6).Formula = "=SUM(D10:E10)"
$ws.Cells.Item(10,7).Value = 3070
$ws.Cells.Item(10,8).Formula = "=SUM(D10:F10)"
$ws.Cells.Item(10,9).Value = 1
$ws.Cells.Item(10,10).Formula = "=(SUM(D10:E10)+5*G10)*I10"

# Clear out the now-unused L10:Q10 cells (keep their existing styles).
$ws.Cells.Item(10,12).ClearContents()
$ws.Cells.Item(10,13).ClearContents()
$ws.Cells.Item(10,14).ClearContents()
$ws.Cells.Item(10,15).ClearContents()
$ws.Cells.Item(10,16).ClearContents()
$ws.Cells.Item(10,17).ClearContents()

# ---------------------------------------------------------------------
# Rows 19-22: add blank placeholder cells in columns A/B to line up
# with the rest of the table (style = default/unstyled, like the A10
# cell above them).
# ---------------------------------------------------------------------
$ws.Cells.Item(19,1).NumberFormat = "GENERAL"
$ws.Cells.Item(19,1).Font.Bold = $false

$ws.Cells.Item(20,1).NumberFormat = "GENERAL"
$ws.Cells.Item(20,1).Font.Bold = $false

$ws.Cells.Item(21,1).NumberFormat = "GENERAL"
$ws.Cells.Item(21,1).Font.Bold = $false
$ws.Cells.Item(21,2).NumberFormat = "GENERAL"
$ws.Cells.Item(21,2).Font.Bold = $false

$ws.Cells.Item(22,1).NumberFormat = "GENERAL"
$ws.Cells.Item(22,1).Font.Bold = $false
$ws.Cells.Item(22,2).NumberFormat = "GENERAL"
$ws.Cells.Item(22,2).Font.Bold = $false

# ---------------------------------------------------------------------
# Row 23 ("MatmultSeq"): was labelled 3x3, now becomes the 2x2 case,
# with updated measurements and a computed (2*2*2) multiplier.
# ---------------------------------------------------------------------
$ws.Cells.Item(23,2).Value = "2x2"
$ws.Cells.Item(23,5).Value = 997
$ws.Cells.Item(23,7).Value = 1956
$ws.Cells.Item(23,9).Formula = "=2*2*2"

# ---------------------------------------------------------------------
# New row 24: the 3x3 case (what row 23 used to represent).
# ---------------------------------------------------------------------
$ws.Cells.Item(24,2).Value = "3x3"
$ws.Cells.Item(24,2).NumberFormat = "GENERAL"
$ws.Cells.Item(24,2).Font.Bold = $true

$ws.Cells.Item(24,4).Value = 1026
$ws.Cells.Item(24,5).Value = 997

$ws.Cells.Item(24,6).NumberFormat = "#,##0"
$ws.Cells.Item(24,6).Font.Bold = $false
$ws.Cells.Item(24,6).Formula = "=SUM(D24:E24)"

$ws.Cells.Item(24,7).Value = 1956

$ws.Cells.Item(24,8).NumberFormat = "#,##0"
$ws.Cells.Item(24,8).Font.Bold = $false
$ws.Cells.Item(24,8).Formula = "=SUM(D24:F24)"

$ws.Cells.Item(24,9).Formula = "=3*3*3"

$ws.Cells.Item(24,10).NumberFormat = "#,##0"
$ws.Cells.Item(24,10).Font.Bold = $true
$ws.Cells.Item(24,10).Formula = "=(SUM(D24:E24)+5*G24)*I24"

# ---------------------------------------------------------------------
# New row 25: the 5x5 case.
# ---------------------------------------------------------------------
$ws.Cells.Item(25,2).Value = "5x5"
$ws.Cells.Item(25,2).NumberFormat = "GENERAL"
$ws.Cells.Item(25,2).Font.Bold = $true

$ws.Cells.Item(25,4).Value = 1026
$ws.Cells.Item(25,5).Value = 997

$ws.Cells.Item(25,6).NumberFormat = "#,##0"
$ws.Cells.Item(25,6).Font.Bold = $false
$ws.Cells.Item(25,6).Formula = "=SUM(D25:E25)"

$ws.Cells.Item(25,7).Value = 1956

$ws.Cells.Item(25,8).NumberFormat = "#,##0"
$ws.Cells.Item(25,8).Font.Bold = $false
$ws.Cells.Item(25,8).Formula = "=SUM(D25:F25)"

$ws.Cells.Item(25,9).Formula = "=5*5*5"

$ws.Cells.Item(25,10).NumberFormat = "#,##0"
$ws.Cells.Item(25,10).Font.Bold = $true
$ws.Cells.Item(25,10).Formula = "=(SUM(D25:E25)+5*G25)*I25"

# ---------------------------------------------------------------------
# Update the view: scroll so column G is the leftmost visible column,
# and select M12 (best-effort - the scroll position is a viewport hint
# that this runtime does not always persist on save).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("M12").Select() | Out-Null
